$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1) - existing columns A-H are unchanged; add new headers
# for the 3 new columns I, J, K.
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 9).Value  = "expected_nohp"
$ws.Cells.Item(1, 10).Value = "expected_alamat"
$ws.Cells.Item(1, 11).Value = "expected_pekerjaan "

# ---------------------------------------------------------------------------
# Row 2 - update address (D2/E2) to "jl kenangan 3"; add I2:K2 = "passed"
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 4).Value = "jl kenangan 3"
$ws.Cells.Item(2, 5).Value = "jl kenangan 3"
$ws.Cells.Item(2, 9).Value  = "passed"
$ws.Cells.Item(2, 10).Value = "passed"
$ws.Cells.Item(2, 11).Value = "passed"

# ---------------------------------------------------------------------------
# Row 3 - update phone number (A3/B3); add I3:K3
# ---------------------------------------------------------------------------
$ws.Cells.Item(3, 1).Value = 81225672116
$ws.Cells.Item(3, 2).Value = 81225672116
$ws.Cells.Item(3, 9).Value  = "failed"
$ws.Cells.Item(3, 10).Value = "passed"
$ws.Cells.Item(3, 11).Value = "passed"

# ---------------------------------------------------------------------------
# Row 4 - update phone number (A4/B4); add I4:K4
# ---------------------------------------------------------------------------
$ws.Cells.Item(4, 1).Value = "081asd"
$ws.Cells.Item(4, 2).Value = 81
$ws.Cells.Item(4, 9).Value  = "failed"
$ws.Cells.Item(4, 10).Value = "passed"
$ws.Cells.Item(4, 11).Value = "passed"

# ---------------------------------------------------------------------------
# Row 5 (new) - full row A-K
# ---------------------------------------------------------------------------
$ws.Cells.Item(5, 1).Value  = 87830815038
$ws.Cells.Item(5, 2).Value  = 87830815038
$ws.Cells.Item(5, 3).Value  = 1998
$ws.Cells.Item(5, 4).Value  = "jl, neraka 7"
$ws.Cells.Item(5, 5).Value  = "jl, neraka 7"
$ws.Cells.Item(5, 6).Value  = "pembantu"
$ws.Cells.Item(5, 7).Value  = "pembantu"
$ws.Cells.Item(5, 8).Value  = "failed"
$ws.Cells.Item(5, 9).Value  = "passed"
$ws.Cells.Item(5, 10).Value = "failed"
$ws.Cells.Item(5, 11).Value = "passed"

# ---------------------------------------------------------------------------
# Row 6 (new) - full row A-K
# ---------------------------------------------------------------------------
$ws.Cells.Item(6, 1).Value  = 87830815038
$ws.Cells.Item(6, 2).Value  = 87830815038
$ws.Cells.Item(6, 3).Value  = 1998
$ws.Cells.Item(6, 4).Value  = "jl kenangan 3"
$ws.Cells.Item(6, 5).Value  = "jl kenangan 3"
$ws.Cells.Item(6, 6).Value  = "dokter?"
$ws.Cells.Item(6, 7).Value  = "dokter"
$ws.Cells.Item(6, 8).Value  = "failed"
$ws.Cells.Item(6, 9).Value  = "passed"
$ws.Cells.Item(6, 10).Value = "passed"
$ws.Cells.Item(6, 11).Value = "failed"

# ---------------------------------------------------------------------------
# Column widths for the new columns I, J, K (best-effort; engine quantizes
# ColumnWidth to 1/6-character steps so an exact match to the source file's
# 1/256-character widths isn't achievable through this property).
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth  = (105/6)
$ws.Columns.Item(10).ColumnWidth = (95/6)
$ws.Columns.Item(11).ColumnWidth = (106/6)

# ---------------------------------------------------------------------------
# View state: scroll so column C is the leftmost visible column, and select
# K7 (the cell just below the newly added data).
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("K7").Select() | Out-Null
